$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 4837.25
$ws.Range("I18").Value = 3599.8
$ws.Range("K18").Value = 3599.8
$ws.Range("M18").Value = -3315.8
$ws.Range("H40").Value = 1625.2084
$ws.Range("I40").Value = 1564.4286
$ws.Range("J40").Value = 1710.3
$ws.Range("K40").Value = 1564.4286
$ws.Range("L40").Value = 1710.3
$ws.Range("M40").Value = -1389.4286
$ws.Range("N40").Value = -2060.3
$ws.Range("H53").Value = 2321.2
$ws.Range("I53").Value = 3258.5715
$ws.Range("J53").Value = 134
$ws.Range("K53").Value = 3258.5715
$ws.Range("L53").Value = 134
$ws.Range("M53").Value = -2621.5715
$ws.Range("N53").Value = -1408
$ws.Range("H55").Value = 380.1
$ws.Range("I55").Value = 357.2857
$ws.Range("J55").Value = 433.33334
$ws.Range("K55").Value = 357.2857
$ws.Range("L55").Value = 433.33334
$ws.Range("M55").Value = -143.2857
$ws.Range("N55").Value = -861.33334
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H69").Value = 4620
$ws.Range("I69").Value = 4620
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 13860
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -12986
$ws.Range("N69").ClearContents()
$ws.Range("H70").Value = 2454.1
$ws.Range("I70").Value = 2285.5715
$ws.Range("J70").Value = 2544.8462
$ws.Range("K70").Value = 6856.7145
$ws.Range("L70").Value = 7634.5386
$ws.Range("M70").Value = -6586.7145
$ws.Range("N70").Value = -8174.5386
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H72").Value = 4620
$ws.Range("I72").Value = 4620
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 41580
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -37212
$ws.Range("N72").ClearContents()
$ws.Range("H73").Value = 2454.1
$ws.Range("I73").Value = 2285.5715
$ws.Range("J73").Value = 2544.8462
$ws.Range("K73").Value = 6856.7145
$ws.Range("L73").Value = 7634.5386
$ws.Range("M73").Value = -5920.7145
$ws.Range("N73").Value = -9506.5386
$ws.Range("H74").Value = 3328.2415
$ws.Range("I74").Value = 3144.2144
$ws.Range("J74").Value = 3500
$ws.Range("K74").Value = 3144.2144
$ws.Range("L74").Value = 3500
$ws.Range("M74").Value = -2208.2144
$ws.Range("N74").Value = -5372
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H76").Value = 6044.4688
$ws.Range("I76").Value = 4641.278
$ws.Range("K76").Value = 4641.278
$ws.Range("M76").Value = -4326.278
$ws.Range("H77").Value = 3328.2415
$ws.Range("I77").Value = 3144.2144
$ws.Range("J77").Value = 3500
$ws.Range("K77").Value = 15721.072
$ws.Range("L77").Value = 17500
$ws.Range("M77").Value = -11041.072
$ws.Range("N77").Value = -26860
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H79").Value = 6044.4688
$ws.Range("I79").Value = 4641.278
$ws.Range("K79").Value = 4641.278
$ws.Range("M79").Value = -3549.278
$ws.Range("H107").Value = 8423.375
$ws.Range("I107").Value = 10458.944
$ws.Range("J107").Value = 2316.6667
$ws.Range("K107").Value = 10458.944
$ws.Range("L107").Value = 2316.6667
$ws.Range("M107").Value = -8538.944
$ws.Range("N107").Value = -6156.6667
$ws.Range("H138").Value = 5684007
$ws.Range("I138").Value = 1795.0625
$ws.Range("K138").Value = 5385.1875
$ws.Range("M138").Value = -245.1875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13864.974
$ws.Range("I32").Value = 14792.151
$ws.Range("K32").Value = 14792.151
$ws.Range("M32").Value = -14505.151
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H122").Value = 11247.8125
$ws.Range("I122").Value = 17488
$ws.Range("J122").Value = 5007.625
$ws.Range("K122").Value = 52464
$ws.Range("L122").Value = 15022.875
$ws.Range("M122").Value = -50014
$ws.Range("N122").Value = -19922.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 12170.444
$ws.Range("I7").Value = 13213
$ws.Range("J7").Value = 3830
$ws.Range("K7").Value = 13213
$ws.Range("L7").Value = 3830
$ws.Range("M7").Value = -13100
$ws.Range("N7").Value = -4056
$ws.Range("H22").Value = 9497.909
$ws.Range("I22").Value = 447.8889
$ws.Range("J22").Value = 50223
$ws.Range("K22").Value = 447.8889
$ws.Range("L22").Value = 50223
$ws.Range("M22").Value = -274.8889
$ws.Range("N22").Value = -50569
$ws.Range("H99").Value = 1000
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 1000
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 1000
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -3996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 4671.6665
$ws.Range("I21").Value = 5000
$ws.Range("J21").Value = 4507.5
$ws.Range("K21").Value = 5000
$ws.Range("L21").Value = 4507.5
$ws.Range("M21").Value = -4765
$ws.Range("N21").Value = -4977.5
$ws.Range("H22").Value = 386.41666
$ws.Range("I22").Value = 124.166664
$ws.Range("J22").Value = 648.6667
$ws.Range("K22").Value = 124.166664
$ws.Range("L22").Value = 648.6667
$ws.Range("M22").Value = 225.833336
$ws.Range("N22").Value = -1348.6667
$ws.Range("H44").Value = 123333.336
$ws.Range("I44").Value = 170000
$ws.Range("K44").Value = 170000
$ws.Range("M44").Value = -169558
$ws.Range("H47").Value = 41402
$ws.Range("I47").Value = 44064
$ws.Range("J47").Value = 40071
$ws.Range("K47").Value = 44064
$ws.Range("L47").Value = 40071
$ws.Range("M47").Value = -43498
$ws.Range("N47").Value = -41203
$ws.Range("H48").Value = 34275.5
$ws.Range("J48").Value = 34275.5
$ws.Range("L48").Value = 34275.5
$ws.Range("N48").Value = -35227.5
$ws.Range("H51").Value = 6049.5
$ws.Range("I51").Value = 1000
$ws.Range("K51").Value = 1000
$ws.Range("M51").Value = -264
$ws.Range("H61").Value = 6049.5
$ws.Range("I61").Value = 1000
$ws.Range("K61").Value = 1000
$ws.Range("M61").Value = -652
$ws.Range("H62").Value = 2383.4783
$ws.Range("I62").Value = 2370.5264
$ws.Range("J62").Value = 2445
$ws.Range("K62").Value = 2370.5264
$ws.Range("L62").Value = 2445
$ws.Range("M62").Value = -1746.5264
$ws.Range("N62").Value = -3693
$ws.Range("H65").Value = 2383.4783
$ws.Range("I65").Value = 2370.5264
$ws.Range("J65").Value = 2445
$ws.Range("K65").Value = 11852.632
$ws.Range("L65").Value = 12225
$ws.Range("M65").Value = -8732.632000000001
$ws.Range("N65").Value = -18465
$ws.Range("H105").Value = 2706.6667
$ws.Range("I105").Value = 1893.3334
$ws.Range("J105").Value = 2977.7778
$ws.Range("K105").Value = 1893.3334
$ws.Range("L105").Value = 2977.7778
$ws.Range("M105").Value = -146.3334
$ws.Range("N105").Value = -6471.7778
$ws.Range("H107").Value = 495
$ws.Range("I107").Value = 417.46155
$ws.Range("J107").Value = 595.8
$ws.Range("K107").Value = 417.46155
$ws.Range("L107").Value = 595.8
$ws.Range("M107").Value = 1502.53845
$ws.Range("N107").Value = -4435.8
$ws.Range("H140").Value = 46900
$ws.Range("J140").Value = 46900
$ws.Range("L140").Value = 46900
$ws.Range("N140").Value = -57260

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 3439
$ws.Range("J80").Value = 3722.2222
$ws.Range("L80").Value = 11166.6666
$ws.Range("N80").Value = -13038.6666
$ws.Range("H83").Value = 3439
$ws.Range("J83").Value = 3722.2222
$ws.Range("L83").Value = 33499.99980000001
$ws.Range("N83").Value = -42859.99980000001
$ws.Range("H118").Value = 7822.75
$ws.Range("J118").Value = 10496
$ws.Range("L118").Value = 31488
$ws.Range("N118").Value = -33974

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2476.276
$ws.Range("I102").Value = 3142
$ws.Range("J102").Value = 1211.4
$ws.Range("K102").Value = 3142
$ws.Range("L102").Value = 1211.4
$ws.Range("M102").Value = -1520
$ws.Range("N102").Value = -4455.4
$ws.Range("H122").Value = 3510753.2
$ws.Range("I122").Value = 4446121
$ws.Range("K122").Value = 13338363
$ws.Range("M122").Value = -13335913
$ws.Range("H126").Value = 3095.5833
$ws.Range("I126").Value = 1888.3182
$ws.Range("K126").Value = 5664.9546
$ws.Range("M126").Value = -3194.9546
$ws.Range("H138").Value = 56516.168
$ws.Range("J138").Value = 56516.168
$ws.Range("L138").Value = 56516.168
$ws.Range("N138").Value = -66796.16800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1509.5
$ws.Range("I46").Value = 1418.1818
$ws.Range("J46").Value = 1621.1111
$ws.Range("K46").Value = 1418.1818
$ws.Range("L46").Value = 1621.1111
$ws.Range("M46").Value = -1230.1818
$ws.Range("N46").Value = -1997.1111
$ws.Range("H132").Value = 9097512
$ws.Range("I132").Value = 3893.8823
$ws.Range("J132").Value = 23820514
$ws.Range("K132").Value = 11681.6469
$ws.Range("L132").Value = 71461542
$ws.Range("M132").Value = -9151.6469
$ws.Range("N132").Value = -71466602

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2524.4736
$ws.Range("I126").Value = 1651.125
$ws.Range("J126").Value = 7182.3335
$ws.Range("K126").Value = 4953.375
$ws.Range("L126").Value = 21547.0005
$ws.Range("M126").Value = -2483.375
$ws.Range("H138").Value = 64494.75
$ws.Range("J138").Value = 64494.75
$ws.Range("L138").Value = 64494.75
$ws.Range("N138").Value = -74774.75
